$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the "_old" / "_new" header suffixes to "_FV2210" / "_FV2304" ---
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
$newCols = @("L","M","N","O","P","Q","R","S","T","U")

foreach ($col in $oldCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = ($cell.Value2 -replace "_old$", "_FV2210")
}

foreach ($col in $newCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = ($cell.Value2 -replace "_new$", "_FV2304")
}

# --- 2. Turn the header + data range into an Excel Table (ListObject) ---
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U76"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# --- 3. Freeze the header row ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
